# Apply USCDI5-Sandbox style metadata update to the "Metadata" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row above the current "Jurisdiction" row (row 11) so a second
# Contact row can be added, pushing Jurisdiction/Description/Purpose/Copyright/
# Immutable down by one row.
$ws.Rows.Item(11).Insert()

# The inserted row doesn't automatically pick up the same cell style as the
# surrounding data rows, so copy formatting down from the row above (row 10,
# the existing Contact row) into the freshly inserted row 11.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update Publisher value (row 9, column B)
$ws.Cells.Item(9, 2).Value = "HL7 International / Cross-Group Projects"

# Update existing Contact row (row 10) to the new HL7 CGP contact details
$ws.Cells.Item(10, 2).Value = "HL7 International / Cross-Group Projects (http://www.hl7.org/Special/committees/cgp, cgp@lists.HL7.org)"

# Populate the newly inserted row 11 with the old Health eData contact info
$ws.Cells.Item(11, 1).Value = "Contact"
$ws.Cells.Item(11, 2).Value = "Health eData Inc (mailto:ehaas@healthedatainc.com)"

# Update Jurisdiction value (now row 12, column B)
$ws.Cells.Item(12, 2).Value = "United States of America"
